$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 160
$ws.Range("C2").Value = 123
$ws.Range("D2").Value = 123
$ws.Range("E2").Value = 112
$ws.Range("F2").Value = 179
$ws.Range("G2").Value = 244
$ws.Range("H2").Value = 160
$ws.Range("I2").Value = 76
$ws.Range("J2").Value = 120

$ws.Range("B3").Value = 28
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 17
$ws.Range("G3").Value = 37
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 19

$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 17
$ws.Range("H5").Value = 11
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 15

$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 11
$ws.Range("H8").Value = 9
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 3

$ws.Range("B9").Value = 6
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 7
$ws.Range("G9").Value = 9
$ws.Range("H9").Value = 9
$ws.Range("I9").Value = 2
$ws.Range("J9").Value = 1

$ws.Range("B10").Value = 34
$ws.Range("C10").Value = 20
$ws.Range("D10").Value = 30
$ws.Range("E10").Value = 24
$ws.Range("F10").Value = 43
$ws.Range("G10").Value = 51
$ws.Range("H10").Value = 25
$ws.Range("I10").Value = 23
$ws.Range("J10").Value = 23

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 14
$ws.Range("G11").Value = 11
$ws.Range("H11").Value = 7
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 6

$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 7
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = 13
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 7

$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 6
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 4

$ws.Range("B14").Value = 8
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1

$ws.Range("B15").Value = 6
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = 9
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 5

$ws.Range("B16").Value = 6
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = 18
$ws.Range("J16").Value = 0

$ws.Range("B17").Value = 23
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 31
$ws.Range("E17").Value = 22
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = 22
$ws.Range("I17").Value = 5
$ws.Range("J17").Value = 23

$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = 0

$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 5
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 3
$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 5

$ws.Range("B21").Value = 9
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 10
$ws.Range("H21").Value = 4
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 15

$ws.Range("B22").Value = 5
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 6
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2

$ws.Range("B23").Value = 3
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 6
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 4
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1

$ws.Range("B24").Value = 20
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 19
$ws.Range("F24").Value = 34
$ws.Range("G24").Value = 29
$ws.Range("H24").Value = 29
$ws.Range("I24").Value = 11
$ws.Range("J24").Value = 21

$ws.Range("B25").Value = 8
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 14
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 8

$ws.Range("B26").Value = 9
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 12
$ws.Range("H26").Value = 10
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 6

$ws.Range("C29").Value = 2
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 1
$ws.Range("J29").Value = 2

$ws.Range("B30").Value = 3
$ws.Range("C30").Value = 6
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 7
$ws.Range("H30").Value = 4
$ws.Range("I30").Value = 6
$ws.Range("J30").Value = 5

$ws.Range("B31").Value = 24
$ws.Range("C31").Value = 25
$ws.Range("D31").Value = 19
$ws.Range("E31").Value = 20
$ws.Range("F31").Value = 33
$ws.Range("G31").Value = 49
$ws.Range("H31").Value = 42
$ws.Range("I31").Value = 19
$ws.Range("J31").Value = 15

$ws.Range("B32").Value = 1
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 12
$ws.Range("F32").Value = 17
$ws.Range("G32").Value = 8
$ws.Range("H32").Value = 8
$ws.Range("I32").Value = 9
$ws.Range("J32").Value = 1

$ws.Range("B33").Value = 2
$ws.Range("C33").Value = 7
$ws.Range("D33").Value = 6
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 5
$ws.Range("G33").Value = 7
$ws.Range("H33").Value = 8
$ws.Range("I33").Value = 4
$ws.Range("J33").Value = 5

$ws.Range("B34").Value = 8
$ws.Range("C34").Value = 5
$ws.Range("D34").Value = 5
$ws.Range("F34").Value = 6
$ws.Range("G34").Value = 12
$ws.Range("H34").Value = 10
$ws.Range("J34").Value = 2

$ws.Range("B35").Value = 8
$ws.Range("C35").Value = 6
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 4
$ws.Range("G35").Value = 14
$ws.Range("H35").Value = 12
$ws.Range("I35").Value = 2
$ws.Range("J35").Value = 2

$ws.Range("B36").Value = 5
$ws.Range("C36").Value = 6
$ws.Range("D36").Value = 5
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 8
$ws.Range("H36").Value = 4
$ws.Range("I36").Value = 4
$ws.Range("J36").Value = 5

$ws.Range("B37").Value = 31
$ws.Range("C37").Value = 29
$ws.Range("D37").Value = 20
$ws.Range("E37").Value = 19
$ws.Range("F37").Value = 28
$ws.Range("G37").Value = 37
$ws.Range("H37").Value = 13
$ws.Range("I37").Value = 13
$ws.Range("J37").Value = 19

$ws.Range("B38").Value = 12
$ws.Range("C38").Value = 12
$ws.Range("D38").Value = 6
$ws.Range("E38").Value = 4
$ws.Range("F38").Value = 15
$ws.Range("G38").Value = 19
$ws.Range("I38").Value = 5
$ws.Range("J38").Value = 4

$ws.Range("B39").Value = 4
$ws.Range("C39").Value = 4
$ws.Range("D39").Value = 4
$ws.Range("E39").Value = 2
$ws.Range("F39").Value = 3
$ws.Range("G39").Value = 7
$ws.Range("H39").Value = 2
$ws.Range("I39").Value = 4
$ws.Range("J39").Value = 4

$ws.Range("B41").Value = 7
$ws.Range("C41").Value = 6
$ws.Range("D41").Value = 4
$ws.Range("E41").Value = 6
$ws.Range("F41").Value = 4
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 4
$ws.Range("J41").Value = 3

$ws.Range("B42").Value = 8
$ws.Range("C42").Value = 7
$ws.Range("D42").Value = 6
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = 6
$ws.Range("G42").Value = 7
$ws.Range("H42").Value = 7
$ws.Range("I42").Value = 4
$ws.Range("J42").Value = 8

